# Quiron 2.0 atualizado (criptografia, cadastro de servidores e tela inicial)
# Adds four new "diario de bordo" entries (rows 28-31) to Planilha1:
#   row 28 (2020-03-12): finish filling the already-started row (C..H)
#   row 29 (2020-03-13), row 30 (2020-03-14): full new rows
#   row 31 (2020-03-15): full new row, only a morning shift (no D/E)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha1")

$xlCenter = -4108

function Set-DateCell($addr, $value) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "m/d/yy"
    $r.Value = $value
}

function Set-TimeCell($addr, $value) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "h:mm"
    $r.Value = $value
}

function Set-TextCell($addr, $value) {
    $r = $ws.Range($addr)
    $r.HorizontalAlignment = $xlCenter
    $r.Value = $value
}

# ---- Row 28 : 12/03/2020 -> finish the row that already had A28/B28 ----
Set-DateCell "A28" 43902
Set-TimeCell "B28" 0.375
Set-TimeCell "C28" 0.45833333333333331
Set-TimeCell "D28" 0.58333333333333337
Set-TimeCell "E28" 0.66666666666666663
$ws.Range("F28").Value = 4
$ws.Range("G28").Formula = "=G27+F28"
Set-TextCell "H28" "Reunião com o professor Leonardo e tentativa de criptografia"

# ---- Row 29 : 13/03/2020 ----
Set-DateCell "A29" 43903
Set-TimeCell "B29" 0.41666666666666669
Set-TimeCell "C29" 0.5
Set-TimeCell "D29" 0.625
Set-TimeCell "E29" 0.75
$ws.Range("F29").Value = 5
$ws.Range("G29").Formula = "=G28+F29"
Set-TextCell "H29" "Conclusão da criptografia e elaboração da página inicial"

# ---- Row 30 : 14/03/2020 ----
Set-DateCell "A30" 43904
Set-TimeCell "B30" 0.41666666666666669
Set-TimeCell "C30" 0.5
Set-TimeCell "D30" 0.625
Set-TimeCell "E30" 0.75
$ws.Range("F30").Value = 5
$ws.Range("G30").Formula = "=G29+F30"
Set-TextCell "H30" "Conclusão da página inicial, ajustes no PDF e no BD e concerto de erros"

# ---- Row 31 : 15/03/2020 (single, morning-only shift) ----
Set-DateCell "A31" 43905
Set-TimeCell "B31" 0.54166666666666663
Set-TimeCell "C31" 0.66666666666666663
$ws.Range("F31").Value = 3
$ws.Range("G31").Value = 126
Set-TextCell "H31" "Elaboração da tela de cadastro de servidores"

# ---- Selection / active cell moves to the newly added last row ----
$ws.Activate()
$ws.Range("A31").Select()
